# Update gh-pages output data (regenerated at 456a3b4)
# Sheet 1 = 展览 (index 1), Sheet 4 = 全部类型 (index 4)

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# 展览 sheet updates
$wsExhibit.Range("F2").Value = 6645
$wsExhibit.Range("F4").Value = 107
$wsExhibit.Range("F7").Value = 80

# 全部类型 sheet updates
$wsAll.Range("F2").Value = 6645
$wsAll.Range("F5").Value = 107
$wsAll.Range("F9").Value = 80
